$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Row 91 F value changes from "1" to "12"; rows 92-111 are brand new order lines
# appended to the bottom of the Orders sheet. Columns B/D/E/G-L stay blank, as
# in the rest of the sheet. Numeric-looking text (PackageID in A, Number in F)
# is forced to Text format first so it round-trips as a string, matching the
# rest of the workbook's convention of storing every value as text.
$rows = @(
    @{ R = 91; A = $null; C = $null; F = '12' },
    @{ R = 92; A = $null; C = '7_翠绿洋桔梗_Dark Green Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '8' },
    @{ R = 93; A = $null; C = '3_波浪白洋桔梗_Wavy White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '10' },
    @{ R = 94; A = $null; C = '12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '5' },
    @{ R = 95; A = $null; C = '9_茶色洋桔梗_Tea Color Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '5' },
    @{ R = 96; A = '13'; C = '110_绣球单瓣浅蓝_Hydrangea Light Blue S_Hydrangea L._1stem'; F = '45' },
    @{ R = 97; A = $null; C = '106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem'; F = '15' },
    @{ R = 98; A = $null; C = '118_绣球老绿_Hydrangea Garden Lace_Hydrangea L._1stem'; F = '65' },
    @{ R = 99; A = $null; C = '322_喷泉草_Grasses Panicum_undefined_1bunch'; F = '5' },
    @{ R = 100; A = $null; C = '479_绿灵草_lepidium_undefined_1bunch'; F = '5' },
    @{ R = 101; A = $null; C = '9_茶色洋桔梗_Tea Color Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '10' },
    @{ R = 102; A = $null; C = '744_永生吊米深红_undefined_undefined_1bunch'; F = '10' },
    @{ R = 103; A = '14'; C = '47_拉丝玫红_Spider Dark Pink_Gerbera L._20stems'; F = '3' },
    @{ R = 104; A = $null; C = '71_霜雪mini_Snowy_Gerbera L._20stems'; F = '5' },
    @{ R = 105; A = $null; C = '45_拉丝艳粉_Spider pink+_Gerbera L._20stems'; F = '2.5' },
    @{ R = 106; A = $null; C = '12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g'; F = '15' },
    @{ R = 107; A = '1'; C = '556_马尾松_undefined_undefined_1bunch'; F = '20' },
    @{ R = 108; A = $null; C = '592_进口春兰叶_undefined_undefined_1bunch'; F = '20' },
    @{ R = 109; A = $null; C = '522_山归来绿_Smilax china_undefined_1bunch'; F = '5' },
    @{ R = 110; A = $null; C = '688_山归来橙_undefined_undefined_1bunch'; F = '5' },
    @{ R = 111; A = $null; C = '439_九星叶_undefined_undefined_1bunch'; F = $null }
)

foreach ($row in $rows) {
    if ($null -ne $row.A) {
        $cell = $ws.Cells.Item($row.R, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row.A
    }
    if ($null -ne $row.C) {
        $ws.Cells.Item($row.R, 3).Value = $row.C
    }
    if ($null -ne $row.F) {
        $cell = $ws.Cells.Item($row.R, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $row.F
    }
}

# Summary sheet: TotalNumber (G2) is a concatenation of every Number value in
# the Orders sheet; it grows as new order lines are appended.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("G2").NumberFormat = "@"
$summary.Range("G2").Value = "055155552510652566555525321515822555510555551255156558101576510612610551051510555510158105151051541156111387865775125551281055451565551010352.5152020550"
